# Applies the "Week 17 interview cards" edit to the Spring Boot proposal doc.
# (Renames a handful of "Resource(r)" bullet items to the "rezource(r)" in-joke
#  spelling, rewords a couple of bullets, and appends the new "Stretch Goal"
#  interview-card bullets for Reviews / amenities / knowledge articles.)
$d = $word.ActiveDocument

# 1) Title line: "Resource Scheduler" -> "Rezource Scheduler (The “z” is on purpose)"
$d.Content.Find.Execute("Resource Scheduler", $true, $false, $false, $false, $false, $true, 1, $false, "Rezource Scheduler (The “z” is on purpose)", 2) | Out-Null

# 2) "Browse all Resources" -> "Register as a rezourcer"
$d.Content.Find.Execute("Browse all Resources", $true, $false, $false, $false, $false, $true, 1, $false, "Register as a rezourcer", 2) | Out-Null

# 3) "Look at Resource Availability" -> "Register as a Scheduler"
$d.Content.Find.Execute("Look at Resource Availability", $true, $false, $false, $false, $false, $true, 1, $false, "Register as a Scheduler", 2) | Out-Null

# 4) "Browse Resources by whether it is a person, place or thing as well as by location." -> "Browse all rezources"
$d.Content.Find.Execute("Browse Resources by whether it is a person, place or thing as well as by location.", $true, $false, $false, $false, $false, $true, 1, $false, "Browse all rezources", 2) | Out-Null

# 5) "Schedule a resource" -> "Look at rezource availability"
$d.Content.Find.Execute("Schedule a resource", $true, $false, $false, $false, $false, $true, 1, $false, "Look at rezource availability", 2) | Out-Null

# 6) "Create a resource" -> "Browse rezources by whether it is a service, place or thing as well as by location."
#    None of the replacements above add/remove paragraphs, so this bullet is
#    still Paragraphs.Item(51) in the original numbering.
$createIdx = 51
$d.Content.Find.Execute("Create a resource", $true, $false, $false, $false, $false, $true, 1, $false, "Browse rezources by whether it is a service, place or thing as well as by location.", 2) | Out-Null

# ...and two brand new bullet siblings right after it: "Schedule a rezource" / "Create a rezource".
$createPara = $d.Paragraphs.Item($createIdx)
$createPara.Range.InsertParagraphAfter()
$scheduleRezIdx = $createIdx + 1
$d.Paragraphs.Item($scheduleRezIdx).Range.Text = "Schedule a rezource"

$d.Paragraphs.Item($scheduleRezIdx).Range.InsertParagraphAfter()
$createRezIdx = $scheduleRezIdx + 1
$d.Paragraphs.Item($createRezIdx).Range.Text = "Create a rezource"

# 7) "Calculate the travel time for a person resource ..." -> "... for a service resource ..."
$d.Content.Find.Execute("Calculate the travel time for a person resource and add it into the Schedule availability", $true, $false, $false, $false, $false, $true, 1, $false, "Calculate the travel time for a service resource and add it into the Schedule availability", 2) | Out-Null

# 8) Append six new "Stretch Goal" bullet paragraphs at the very end of the document:
#    three ilvl-0 bullets (Reviews / amenities / knowledge articles) and three
#    ilvl-1 sub-bullets underneath "knowledge articles".
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Range.ListFormat.ListLevelNumber = 1
$d.Paragraphs.Item($idx).Range.Text = "Allow people to leave Reviews for a rezource"

$d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Range.ListFormat.ListLevelNumber = 1
$d.Paragraphs.Item($idx).Range.Text = "Have an amenities table/list for a rezource"

$d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Range.ListFormat.ListLevelNumber = 1
$d.Paragraphs.Item($idx).Range.Text = "Have “knowledge articles”"

$d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs.Item($idx).Range.Text = "Things to know"

$d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs.Item($idx).Range.Text = "Safety considerations"

$d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
$idx = $d.Paragraphs.Count
$d.Paragraphs.Item($idx).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs.Item($idx).Range.Text = "Cancellation policy"

Write-Output "Edit complete."
